$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.867.90'
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").Value = '2.317.75'
$ws.Range("E3").Value = '  +1.47%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("E7").Value = '  +0.33%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.80%  '

$ws.Range("E12").Value = '  +0.48%  '

$ws.Range("E13").Value = '  +0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("D15").Value = '2.681.96'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").Value = '2.320.77'
$ws.Range("E16").Value = '  +2.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.792'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.54%  '

$ws.Range("D18").Value = '42.801.82'
$ws.Range("E18").Value = '  +0.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.33%  '

$ws.Range("E20").Value = '  +3.26%  '

$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.92'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.16%  '

$ws.Range("B24").Value = 'ImmutableX'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.06%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +0.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.35'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -1.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.13%  '

$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.09%  '

$ws.Range("E34").Value = '  +0.17%  '

$ws.Range("E35").Value = '  +0.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0700'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.50%  '

$ws.Range("E37").Value = '  -1.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.95%  '

$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.87'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.81%  '

$ws.Range("D43").Value = '1.934.58'
$ws.Range("E43").Value = '  -2.87%  '

$ws.Range("E44").Value = '  +0.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '

$ws.Range("E46").Value = '  +2.69%  '

$ws.Range("E47").Value = '  +0.42%  '

$ws.Range("D48").Value = '2.549.23'
$ws.Range("E48").Value = '  +1.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("E50").Value = '  -3.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.07'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.55%  '
